$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66
$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Den Haag Bezuidenhoutseweg"
$ws.Cells.Item($row, 3).Value = "KDV"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-08-29"
$ws.Cells.Item($row, 4).ClearFormats()
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
